$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Replace the data table (rows 2-17). The old "Account / Rep / Manager /
#    Product / Quantity / Price / Status" rows are swapped for a new set of
#    rows: new account numbers in A, a survey-question label in B, the rep
#    name moves to D, a 1-5 score in F and a rating label in H. Columns
#    C, E and G are no longer used in the new data.
# ---------------------------------------------------------------------------
$rows = @(
  @(2,  105209, "Week Password",                 "Juan Para",     3, "Good"),
  @(3,  105209, "CBTS overdue",                   "Juan Para",     3, "Good"),
  @(4,  105209, "PhishMe clicked",                "Juan Para",     2, "Bad"),
  @(5,  105209, "Security Incident Involvement",  "Juan Para",     1, "Very bad"),
  @(6,  104822, "Week Password",                  "Felipe Fiorin", 5, "Very good"),
  @(7,  104822, "CBTS overdue",                    "Felipe Fiorin", 2, "Bad"),
  @(8,  104822, "PhishMe clicked",                "Felipe Fiorin", 4, "Good"),
  @(9,  104822, "Security Incident Involvement",  "Felipe Fiorin", 2, "Bad"),
  @(10, 978699, "Week Password",                  "Kevin Whelan",  1, "Very bad"),
  @(11, 978699, "CBTS overdue",                    "Kevin Whelan",  5, "Very good"),
  @(12, 978699, "PhishMe clicked",                "Kevin Whelan",  4, "Good"),
  @(13, 978699, "Security Incident Involvement",  "Kevin Whelan",  1, "Bad"),
  @(14, 958039, "Week Password",                  "Maeve Morris",  4, "Good"),
  @(15, 958039, "CBTS overdue",                    "Maeve Morris",  5, "Very good"),
  @(16, 958039, "PhishMe clicked",                "Maeve Morris",  5, "Very good"),
  @(17, 958039, "Security Incident Involvement",  "Maeve Morris",  1, "Very bad")
)

foreach ($r in $rows) {
  $rownum = $r[0]
  $ws.Cells.Item($rownum, 1).Value = $r[1]   # A: Account
  $ws.Cells.Item($rownum, 2).Value = $r[2]   # B: question / name
  $ws.Cells.Item($rownum, 3).ClearContents() # C: no longer used
  $ws.Cells.Item($rownum, 4).Value = $r[3]   # D: rep name
  $ws.Cells.Item($rownum, 5).ClearContents() # E: no longer used
  $ws.Cells.Item($rownum, 6).Value = $r[4]   # F: score
  $ws.Cells.Item($rownum, 7).ClearContents() # G: no longer used
  $ws.Cells.Item($rownum, 8).Value = $r[5]   # H: status / rating
}

# ---------------------------------------------------------------------------
# 2) Page setup is now explicit: paper size 9 (A4), portrait orientation.
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# ---------------------------------------------------------------------------
# 3) Selection / active cell moved to A18 (just past the last data row).
# ---------------------------------------------------------------------------
$ws.Range("A18").Select()

# ---------------------------------------------------------------------------
# 4) The workbook now also carries the (unused) built-in Hyperlink /
#    Followed-Hyperlink cell styles. Apply them briefly to a scratch cell so
#    the style definitions get created, then discard the scratch rows again
#    so the visible sheet data is untouched.
# ---------------------------------------------------------------------------
$ws.Range("Z100").Style = "Hyperlink"
$ws.Range("Z101").Style = "Followed Hyperlink"
$ws.Rows.Item(101).Delete()
$ws.Rows.Item(100).Delete()

Write-Output "edit applied"
